$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.136810279667423
$ws.Range("C2").Value = 0.00226757369614512
$ws.Range("D2").Value = 0.00377928949357521
$ws.Range("E2").Value = 0.00151171579743008
$ws.Range("F2").Value = 0.965230536659108
$ws.Range("G2").Value = 0.00680272108843537
$ws.Range("H2").Value = 0.974300831443689
$ws.Range("I2").Value = 0.00151171579743008
$ws.Range("J2").Value = 0.00529100529100529
$ws.Range("K2").Value = 0.977324263038549
$ws.Range("L2").Value = 0.00604686318972033
$ws.Range("M2").Value = 0.99244142101285
$ws.Range("N2").Value = 0.0158730158730159
$ws.Range("O2").Value = 0.00377928949357521
$ws.Range("P2").Value = 0.054421768707483
$ws.Range("Q2").Value = 0.00377928949357521
$ws.Range("R2").Value = 0.983371126228269
$ws.Range("S2").Value = 0.00453514739229025
$ws.Range("T2").Value = 0.0143613000755858
$ws.Range("U2").Value = 0.00226757369614512
$ws.Range("V2").Value = 0.00226757369614512
$ws.Range("W2").Value = 0.0105820105820106
$ws.Range("X2").Value = 0.00680272108843537

$ws.Range("B3").Value = 0.00377928949357521
$ws.Range("C3").Value = 0.000755857898715042
$ws.Range("D3").Value = 0.0256991685563114
$ws.Range("E3").Value = 0.00453514739229025
$ws.Range("F3").Value = 0.000755857898715042
$ws.Range("G3").Value = 0.0136054421768707
$ws.Range("H3").Value = 0.00226757369614512
$ws.Range("I3").Value = 0.993197278911565
$ws.Range("J3").Value = 0.990929705215419
$ws.Range("K3").Value = 0.018896447467876
$ws.Range("L3").Value = 0.00680272108843537
$ws.Range("M3").Value = 0.00529100529100529
$ws.Range("N3").Value = 0.00226757369614512
$ws.Range("O3").Value = 0.866969009826153
$ws.Range("P3").Value = 0.0513983371126228
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0.00680272108843537
$ws.Range("T3").Value = 0.984126984126984
$ws.Range("U3").Value = 0.994708994708995
$ws.Range("V3").Value = 0.0513983371126228
$ws.Range("W3").Value = 0.0143613000755858
$ws.Range("X3").Value = 0.00377928949357521

$ws.Range("B4").Value = 0.762660619803477
$ws.Range("C4").Value = 0.99244142101285
$ws.Range("D4").Value = 0.0090702947845805
$ws.Range("E4").Value = 0.0204081632653061
$ws.Range("F4").Value = 0.0294784580498866
$ws.Range("G4").Value = 0.973544973544973
$ws.Range("H4").Value = 0.0211640211640212
$ws.Range("I4").Value = 0.000755857898715042
$ws.Range("J4").Value = 0.00226757369614512
$ws.Range("K4").Value = 0.000755857898715042
$ws.Range("L4").Value = 0.983371126228269
$ws.Range("M4").Value = 0.00226757369614512
$ws.Range("N4").Value = 0.000755857898715042
$ws.Range("O4").Value = 0.000755857898715042
$ws.Range("P4").Value = 0.00377928949357521
$ws.Range("Q4").Value = 0.996220710506425
$ws.Range("R4").Value = 0.0158730158730159
$ws.Range("S4").Value = 0.987150415721844
$ws.Range("T4").Value = 0.000755857898715042
$ws.Range("U4").Value = 0.00151171579743008
$ws.Range("V4").Value = 0.0136054421768707
$ws.Range("W4").Value = 0.973544973544973
$ws.Range("X4").Value = 0.987150415721844

$ws.Range("B5").Value = 0.0959939531368103
$ws.Range("C5").Value = 0.00377928949357521
$ws.Range("D5").Value = 0.961451247165533
$ws.Range("E5").Value = 0.973544973544973
$ws.Range("F5").Value = 0.00453514739229025
$ws.Range("G5").Value = 0.00604686318972033
$ws.Range("H5").Value = 0.00226757369614512
$ws.Range("I5").Value = 0.00453514739229025
$ws.Range("J5").Value = 0.00151171579743008
$ws.Range("K5").Value = 0.00302343159486017
$ws.Range("L5").Value = 0.00377928949357521
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0.981103552532124
$ws.Range("O5").Value = 0.128495842781557
$ws.Range("P5").Value = 0.890400604686319
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0.000755857898715042
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0.000755857898715042
$ws.Range("U5").Value = 0.00151171579743008
$ws.Range("V5").Value = 0.932728647014361
$ws.Range("W5").Value = 0.000755857898715042
$ws.Range("X5").Value = 0.00151171579743008
